$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit message: "Added Row 6 to the excel"
# Append a new data row (sheet row 7) to the employee table:
#   No = 6, Name = Minal, Job = RPA Consultant
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Minal"
$ws.Range("C7").Value = "RPA Consultant"

# Move/save the selection as it was left after the edit
$ws.Range("C22").Select()
